# Adjust move balance, new camera logic, modify combo system, add combo damage scaling
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Move balance adjustments (top summary table, rows 3-4) ---
# F3 ( slp impact value) 29 -> 25  (C3 recalculates 12 -> 8 via formula)
$ws.Range("F3").Value = 25
# F4 ( srp impact value) 27 -> 33  (C4 recalculates 4 -> 10 via formula)
$ws.Range("F4").Value = 33

# --- Combo system / damage scaling adjustments (detail table, rows 23-28) ---
# C23 13 -> 8   (F23 recalculates 30 -> 25 via formula)
$ws.Range("C23").Value = 8
# C24 4 -> 10   (F24 recalculates 27 -> 33 via formula)
$ws.Range("C24").Value = 10
# C26 12 -> 14  (F26 recalculates 33 -> 35 via formula)
$ws.Range("C26").Value = 14
# C28 14 -> 15  (F28 recalculates 40 -> 41 via formula)
$ws.Range("C28").Value = 15

# --- New camera logic: update the visible scroll position / active selection ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select()
